# mockup e diagrama att
# Insert a new blank slide at position 3 (between the current 2nd and 3rd
# slides), pushing the former 3rd/4th slides down to positions 4/5.

$p = $ppt.ActivePresentation

# ppLayoutBlank = 12
$s = $p.Slides.Add(3, 12)
